$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 53) down to the
# four new rows (54-57) so the new cells pick up the same style (s="1",
# centered horizontal/vertical alignment) without creating new/orphan
# style entries in styles.xml.
$srcFormat = $ws.Range("A53:E53")
$dstFormat = $ws.Range("A54:E57")
$srcFormat.Copy()
$dstFormat.PasteSpecial(-4122)

# New game (game_id 15): Warhammer game with four players and their decks/positions.
$newRows = @(
    @(15, "Alex",   "Warhammer", "Imperium",   4),
    @(15, "Chris",  "Warhammer", "Chaos",      3),
    @(15, "Kevin",  "Warhammer", "Tyrannids",  2),
    @(15, "Sandro", "Warhammer", "Necron",     1)
)

$startRow = 54
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}

# Match the updated selection recorded in the author's save: the newly
# entered last data cell (E56) is left selected.
$ws.Range("E56").Select()
